# "Guid no nome do arquivo e lista de resultado de avaliacao"
#
# Insert a new column (results of a second evaluation round - "Aquela
# Parada 2" / "&&&&") right after the existing "Aquela Parada" column,
# pushing the old "Prime/Pincel" and "Idioma" columns one slot to the
# right, and point the selection at the freshly filled column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before column D. This shifts the old D ("Prime"/
# "Pincel") and E ("Idioma") columns to E/F respectively, and also moves
# the styled-but-empty cell that used to be F8 to G8 automatically.
$ws.Columns("D:D").Insert()

# Populate the new column with the second evaluation's results.
$ws.Range("D1").Value = "Aquela Parada 2"
$ws.Range("D2").Value = "&&&&"
$ws.Range("D3").Value = "&&&&"
$ws.Range("D4").Value = "&&&&"

# Give the new column a wider, manually sized (non bestFit) width.
$ws.Columns("D:D").ColumnWidth = 24.3

# Leave the selection on the newly added column.
$ws.Range("D1:D4").Select() | Out-Null
